# Drop in results from RMI script
$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$locu  = $wb.Worksheets.Item("OCCF-DpLOCU")
$mocu  = $wb.Worksheets.Item("OCCF-DpMOCU")
$socu  = $wb.Worksheets.Item("OCCF-DpSOCU")

# Update the "2018 dollars" -> "2019 dollars" text labels on the About sheet,
# in the same order the RMI script wrote them.
$about.Range("B26").Value = "2019 dollars per 2012 dollar"
$about.Range("B29").Value = 'which in this case is "2012 dollars per 2019 dollar."'
$about.Range("A21").Value = "million 2019 dollars"
$about.Range("A18").Value = "billion 2019 dollars"
$about.Range("A24").Value = "2018 dollars"

# New conversion factor from the RMI script.
$about.Range("A26").Value = 0.89805481563188172

# Recalculate dependent sheets.
$locu.Range("B2").Formula = '=10^9*About!$A$26'
$mocu.Range("B2").Formula = '=10^6*About!$A$26'
$socu.Range("B2").Formula = '=1*About!A26'

# Turn on iterative calculation, matching the refreshed workbook settings.
$excel.Iteration = $true
$excel.MaxChange = 0.00001

# Leave the cursor where the RMI script left it.
$about.Activate() | Out-Null
$about.Range("A19").Select() | Out-Null

$wb.Save()
